$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9118026264606556
$ws.Range("C2").Value = 0.1825483712090374
$ws.Range("D2").Value = 0.2640623037707428
$ws.Range("F2").Value = 1.277823576529613
$ws.Range("G2").Value = 0.002439361267966054
$ws.Range("J2").Value = 0.2658929057621009
$ws.Range("M2").Value = 0.4265380062410316
$ws.Range("O2").Value = 2.812373875057034

$ws.Range("B3").Value = 0.8118824727546325
$ws.Range("C3").Value = 0.1594099981672343
$ws.Range("D3").Value = 0.2610606059362368
$ws.Range("F3").Value = 1.283391422492976
$ws.Range("G3").Value = 0.002442404532726988
$ws.Range("J3").Value = 0.2670119057478786
$ws.Range("M3").Value = 0.3975929767683013
$ws.Range("O3").Value = 2.833792051614495

$ws.Range("B4").Value = 0.7504831558867409
$ws.Range("C4").Value = 0.1451453340692126
$ws.Range("D4").Value = 0.2593048297493112
$ws.Range("F4").Value = 1.287724246879556
$ws.Range("G4").Value = 0.002444372487571464
$ws.Range("J4").Value = 0.2678553944594455
$ws.Range("M4").Value = 0.379910746431058
$ws.Range("O4").Value = 2.849204723590987

$ws.Range("B5").Value = 0.7254517000858129
$ws.Range("C5").Value = 0.1393182331844685
$ws.Range("D5").Value = 0.2586113639064678
$ws.Range("F5").Value = 1.289719476487335
$ws.Range("G5").Value = 0.002445199510626944
$ws.Range("J5").Value = 0.2682384234942603
$ws.Range("M5").Value = 0.3727281199057515
$ws.Range("O5").Value = 2.856053146210868

$ws.Range("B6").Value = 0.721294636666812
$ws.Range("C6").Value = 0.1383498042830524
$ws.Range("D6").Value = 0.2584975471957307
$ws.Range("F6").Value = 1.290064639410971
$ws.Range("G6").Value = 0.002445338353310217
$ws.Range("J6").Value = 0.2683043978942266
$ws.Range("M6").Value = 0.3715368516618298
$ws.Range("O6").Value = 2.85722457238252

$ws.Range("B7").Value = 0.7501456146275984
$ws.Range("C7").Value = 0.1450668045085592
$ws.Range("D7").Value = 0.259295388131406
$ws.Range("F7").Value = 1.287750226132054
$ws.Range("G7").Value = 0.002444383539439702
$ws.Range("J7").Value = 0.2678604010389876
$ws.Range("M7").Value = 0.3798137853051173
$ws.Range("O7").Value = 2.84929478698345

$ws.Range("B8").Value = 0.8773608630789909
$ws.Range("C8").Value = 0.1745824199531683
$ws.Range("D8").Value = 0.2630092503694499
$ws.Range("F8").Value = 1.279553456031152
$ws.Range("G8").Value = 0.002440390005604114
$ws.Range("J8").Value = 0.2662462526084894
$ws.Range("M8").Value = 0.4165392584951491
$ws.Range("O8").Value = 2.819288743898426

$ws.Range("B9").Value = 1.126404091442225
$ws.Range("C9").Value = 0.2319927744987069
$ws.Range("D9").Value = 0.2709817036192419
$ws.Range("F9").Value = 1.27074777120886
$ws.Range("G9").Value = 0.002433343742510701
$ws.Range("J9").Value = 0.2643236935694873
$ws.Range("M9").Value = 0.4892610744984154
$ws.Range("O9").Value = 2.778445948186118

$ws.Range("B10").Value = 1.309075033723786
$ws.Range("C10").Value = 0.2738731922145234
$ws.Range("D10").Value = 0.27725638636268
$ws.Range("F10").Value = 1.268731468930682
$ws.Range("G10").Value = 0.002428640558682088
$ws.Range("J10").Value = 0.2636714649133438
$ws.Range("M10").Value = 0.5431083509256922
$ws.Range("O10").Value = 2.759486603809592

$ws.Range("B11").Value = 1.392103937271202
$ws.Range("C11").Value = 0.2928584075174001
$ws.Range("D11").Value = 0.2802009335910469
$ws.Range("F11").Value = 1.268786108348294
$ws.Range("G11").Value = 0.002426602780314649
$ws.Range("J11").Value = 0.263540430002621
$ws.Range("M11").Value = 0.5676940190936222
$ws.Range("O11").Value = 2.753276781755119

$ws.Range("B12").Value = 1.423533879865545
$ws.Range("C12").Value = 0.3000377840575084
$ws.Range("D12").Value = 0.2813288560656986
$ws.Range("F12").Value = 1.268946919616312
$ws.Range("G12").Value = 0.00242584567575551
$ws.Range("J12").Value = 0.2635146766623109
$ws.Range("M12").Value = 0.577016685105292
$ws.Range("O12").Value = 2.751273823175751

$ws.Range("B13").Value = 1.416765403730267
$ws.Range("C13").Value = 0.2984920234803212
$ws.Range("D13").Value = 0.2810853659071171
$ws.Range("F13").Value = 1.268906048294298
$ws.Range("G13").Value = 0.002426008085481041
$ws.Range("J13").Value = 0.2635191609237566
$ws.Range("M13").Value = 0.5750083273854187
$ws.Range("O13").Value = 2.751689672742344

$ws.Range("B14").Value = 1.394689931880237
$ws.Range("C14").Value = 0.2934492598795941
$ws.Range("D14").Value = 0.2802934706405438
$ws.Range("F14").Value = 1.268796528063518
$ws.Range("G14").Value = 0.002426540201544759
$ws.Range("J14").Value = 0.2635378326823528
$ws.Range("M14").Value = 0.5684607499008223
$ws.Range("O14").Value = 2.753105003870473

$ws.Range("B15").Value = 1.38116655565193
$ws.Range("C15").Value = 0.2903591199290076
$ws.Range("D15").Value = 0.2798100877738676
$ws.Range("F15").Value = 1.268747702473476
$ws.Range("G15").Value = 0.002426868031593211
$ws.Range("J15").Value = 0.2635523791351062
$ws.Range("M15").Value = 0.5644518004147301
$ws.Range("O15").Value = 2.754017367483328

$ws.Range("B16").Value = 1.303647392708285
$ws.Range("C16").Value = 0.272631096390711
$ws.Range("D16").Value = 0.2770657608175924
$ws.Range("F16").Value = 1.268747485966372
$ws.Range("G16").Value = 0.002428775773135934
$ws.Range("J16").Value = 0.2636833657586379
$ws.Range("M16").Value = 0.5415034041367477
$ws.Range("O16").Value = 2.759941135843206

$ws.Range("B17").Value = 1.256073234429437
$ws.Range("C17").Value = 0.261738261370823
$ws.Range("D17").Value = 0.2754052415799606
$ws.Range("F17").Value = 1.268996532972068
$ws.Range("G17").Value = 0.002429972114885802
$ws.Range("J17").Value = 0.2638061844237001
$ws.Range("M17").Value = 0.5274481598037681
$ws.Range("O17").Value = 2.764194617111059

$ws.Range("B18").Value = 1.228703421188015
$ws.Range("C18").Value = 0.2554667626954199
$ws.Range("D18").Value = 0.2744586460452751
$ws.Range("F18").Value = 1.269231232887449
$ws.Range("G18").Value = 0.0024306697977826
$ws.Range("J18").Value = 0.2638924172143007
$ws.Range("M18").Value = 0.5193724739441734
$ws.Range("O18").Value = 2.766868328740912

$ws.Range("B19").Value = 1.219435412980715
$ws.Range("C19").Value = 0.2533422825132732
$ws.Range("D19").Value = 0.2741396058496974
$ws.Range("F19").Value = 1.269326394465523
$ws.Range("G19").Value = 0.002430907668676546
$ws.Range("J19").Value = 0.2639242905418868
$ws.Range("M19").Value = 0.5166396600660619
$ws.Range("O19").Value = 2.76781258716349

$ws.Range("B20").Value = 1.261138258810377
$ws.Range("C20").Value = 0.2628984694943881
$ws.Range("D20").Value = 0.2755811284527283
$ws.Range("F20").Value = 1.268960553610242
$ws.Range("G20").Value = 0.002429843771367587
$ws.Range("J20").Value = 0.2637914962652061
$ws.Range("M20").Value = 0.5289434856022552
$ws.Range("O20").Value = 2.763718300667989

$ws.Range("B21").Value = 1.401174350675717
$ws.Range("C21").Value = 0.2949307136918833
$ws.Range("D21").Value = 0.280525720370278
$ws.Range("F21").Value = 1.268824891056255
$ws.Range("G21").Value = 0.002426383511507818
$ws.Range("J21").Value = 0.2635317002469009
$ws.Range("M21").Value = 0.5703835929322736
$ws.Range("O21").Value = 2.752679816399336

$ws.Range("B22").Value = 1.492629532723072
$ws.Range("C22").Value = 0.3158076936605028
$ws.Range("D22").Value = 0.2838323722778142
$ws.Range("F22").Value = 1.269553127110854
$ws.Range("G22").Value = 0.002424206854580926
$ws.Range("J22").Value = 0.2635010426296844
$ws.Range("M22").Value = 0.5975404114990113
$ws.Range("O22").Value = 2.74749784535652

$ws.Range("B23").Value = 1.4438247159257
$ws.Range("C23").Value = 0.3046706776428323
$ws.Range("D23").Value = 0.2820607064616212
$ws.Range("F23").Value = 1.269089586328761
$ws.Range("G23").Value = 0.002425360839019793
$ws.Range("J23").Value = 0.263504659833167
$ws.Range("M23").Value = 0.5830397207279958
$ws.Range("O23").Value = 2.750077171821687

$ws.Range("B24").Value = 1.258848419734704
$ws.Range("C24").Value = 0.2623739676064929
$ws.Range("D24").Value = 0.2755015848895965
$ws.Range("F24").Value = 1.268976534825796
$ws.Range("G24").Value = 0.002429901764605095
$ws.Range("J24").Value = 0.2637980881218382
$ws.Range("M24").Value = 0.5282674336401669
$ws.Range("O24").Value = 2.763932932187885

$ws.Range("B25").Value = 1.059081227919478
$ws.Range("C25").Value = 0.2165133599256706
$ws.Range("D25").Value = 0.2687514385570751
$ws.Range("F25").Value = 1.272349520224026
$ws.Range("G25").Value = 0.002435166405118721
$ws.Range("J25").Value = 0.264710475666476
$ws.Range("M25").Value = 0.4695136739142214
$ws.Range("O25").Value = 2.787560308893461

